$wb = $excel.ActiveWorkbook

# "select" is currently the last / active sheet in the workbook.
$select = $wb.Worksheets.Item("select")

# ---------------------------------------------------------------------
# Add the new "tab" worksheet right after "select".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $select)
$newSheet.Name = "tab"

# ---------------------------------------------------------------------
# Copy cell formatting (styles) from existing cells on "select" so the
# new sheet's cells reuse the same style indices instead of new ones.
# (-4122 = xlPasteFormats). Only paste onto cells that must actually end
# up with a *non-default* style, one cell at a time, so we never create
# phantom empty cells on rows/cols that should stay untouched.
# ---------------------------------------------------------------------

# style 3 -> header row, columns A:D
$select.Range("A1").Copy()
[void]$newSheet.Range("A1:D1").PasteSpecial(-4122)

# style 4 -> header row, column E
$select.Range("E1").Copy()
[void]$newSheet.Range("E1").PasteSpecial(-4122)

# style 6
$select.Range("B3").Copy()
[void]$newSheet.Range("C3").PasteSpecial(-4122)
[void]$newSheet.Range("D3").PasteSpecial(-4122)
[void]$newSheet.Range("D4").PasteSpecial(-4122)
[void]$newSheet.Range("D12").PasteSpecial(-4122)

# style 7
$select.Range("E5").Copy()
[void]$newSheet.Range("E4").PasteSpecial(-4122)
[void]$newSheet.Range("E12").PasteSpecial(-4122)

# style 2 (blank, formatted cell)
$select.Range("E2").Copy()
[void]$newSheet.Range("E3").PasteSpecial(-4122)

# style 1 (blank, formatted cell)
$select.Range("B2").Copy()
[void]$newSheet.Range("C4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------
$newSheet.Range("A1").Value = "Desc"
$newSheet.Range("B1").Value = "Steps"
$newSheet.Range("C1").Value = "Locator"
$newSheet.Range("D1").Value = "Action"
$newSheet.Range("E1").Value = "Data"

# ---------------------------------------------------------------------
# Test-case body ("sheet – action – tab")
# ---------------------------------------------------------------------
$newSheet.Range("A2").Value = "sheet – action – tab"

$newSheet.Range("C3").Value = "https://letcode.in/windows"
$newSheet.Range("D3").Value = "url"

$newSheet.Range("D4").Value = "title"
$newSheet.Range("E4").Value = "Window"

$newSheet.Range("C5").Value = "button#home"
$newSheet.Range("D5").Value = "click:tab"

$newSheet.Range("D6").Value = "title"
$newSheet.Range("E6").Value = "Testing Hub"

$newSheet.Range("C7").Value = "//a[@title='Koushik Chatterjee']"
$newSheet.Range("D7").Value = "click:tab"

$newSheet.Range("D8").Value = "title"
$newSheet.Range("E8").Value = "LinkedIn"

$newSheet.Range("D9").Value = "tab:back"

$newSheet.Range("D10").Value = "title"
$newSheet.Range("E10").Value = "Testing Hub"

$newSheet.Range("D11").Value = "tab:back"

$newSheet.Range("D12").Value = "title?"
$newSheet.Range("E12").Value = "Window"

# ---------------------------------------------------------------------
# Column widths, approximating the "select" sheet's layout
# (C ~35.06, D ~11.02, E ~17.09 OOXML character-width units).
# ---------------------------------------------------------------------
$newSheet.Range("C:C").ColumnWidth = 34.166666666666664
$newSheet.Range("D:D").ColumnWidth = 10.166666666666666
$newSheet.Range("E:E").ColumnWidth = 16.333333333333332

# Row 5 on the new sheet uses the slightly shorter "12.8" row height.
$newSheet.Rows.Item(5).RowHeight = 12.8

# ---------------------------------------------------------------------
# Update the "select" sheet's remembered selection (cursor moved to C4,
# it is no longer the active tab).
# ---------------------------------------------------------------------
[void]$select.Range("C4").Select()

# ---------------------------------------------------------------------
# Finally, make the new "tab" sheet the active sheet/tab with its
# cursor at J18 (this must run last so "tab" ends up as the active tab).
# ---------------------------------------------------------------------
[void]$newSheet.Range("J18").Select()
